$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item(1)

$ws.Range("B2").Value = 0.5666666666666667
$ws.Range("C2").Value = 0.8947368421052632
$ws.Range("D2").Value = 0.6938775510204083
$ws.Range("B3").Value = 0.8787878787878788
$ws.Range("C3").Value = 0.5272727272727272
$ws.Range("D3").Value = 0.6590909090909091
$ws.Range("B4").Value = 0.6774193548387096
$ws.Range("C4").Value = 0.6774193548387096
$ws.Range("D4").Value = 0.6774193548387096
$ws.Range("E4").Value = 0.6774193548387096
$ws.Range("B5").Value = 0.7227272727272727
$ws.Range("C5").Value = 0.7110047846889952
$ws.Range("D5").Value = 0.6764842300556586
$ws.Range("B6").Value = 0.7512544802867385
$ws.Range("C6").Value = 0.6774193548387096
$ws.Range("D6").Value = 0.6733048057932851
$ws.Range("B7").Value = 0.6923076923076923
$ws.Range("C7").Value = 0.4736842105263158
$ws.Range("D7").Value = 0.5625
$ws.Range("B8").Value = 0.7014925373134329
$ws.Range("C8").Value = 0.8545454545454545
$ws.Range("D8").Value = 0.7704918032786885
$ws.Range("B9").Value = 0.6989247311827957
$ws.Range("C9").Value = 0.6989247311827957
$ws.Range("D9").Value = 0.6989247311827957
$ws.Range("E9").Value = 0.6989247311827957
$ws.Range("B10").Value = 0.6969001148105626
$ws.Range("C10").Value = 0.6641148325358851
$ws.Range("D10").Value = 0.6664959016393442
$ws.Range("B11").Value = 0.6977395898917323
$ws.Range("C11").Value = 0.6989247311827957
$ws.Range("D11").Value = 0.6855059051648158
$ws.Range("B12").Value = 0.7142857142857143
$ws.Range("C12").Value = 0.131578947368421
$ws.Range("D12").Value = 0.2222222222222222
$ws.Range("B13").Value = 0.6162790697674418
$ws.Range("C13").Value = 0.9636363636363636
$ws.Range("D13").Value = 0.75177304964539
$ws.Range("B14").Value = 0.6236559139784946
$ws.Range("C14").Value = 0.6236559139784946
$ws.Range("D14").Value = 0.6236559139784946
$ws.Range("E14").Value = 0.6236559139784946
$ws.Range("B15").Value = 0.665282392026578
$ws.Range("C15").Value = 0.5476076555023923
$ws.Range("D15").Value = 0.4869976359338061
$ws.Range("B16").Value = 0.6563247954845854
$ws.Range("C16").Value = 0.6236559139784946
$ws.Range("D16").Value = 0.5353974427413
$ws.Range("B17").Value = 0.6666666666666666
$ws.Range("C17").Value = 0.7368421052631579
$ws.Range("D17").Value = 0.7
$ws.Range("B18").Value = 0.803921568627451
$ws.Range("C18").Value = 0.7454545454545455
$ws.Range("D18").Value = 0.7735849056603775
$ws.Range("B19").Value = 0.7419354838709677
$ws.Range("C19").Value = 0.7419354838709677
$ws.Range("D19").Value = 0.7419354838709677
$ws.Range("E19").Value = 0.7419354838709677
$ws.Range("B20").Value = 0.7352941176470589
$ws.Range("C20").Value = 0.7411483253588517
$ws.Range("D20").Value = 0.7367924528301888
$ws.Range("B21").Value = 0.7478389205144423
$ws.Range("C21").Value = 0.7419354838709677
$ws.Range("D21").Value = 0.7435179549604384
$ws.Range("B22").Value = 0.7555555555555555
$ws.Range("C22").Value = 0.8947368421052632
$ws.Range("D22").Value = 0.8192771084337349
$ws.Range("B23").Value = 0.9166666666666666
$ws.Range("C23").Value = 0.8
$ws.Range("D23").Value = 0.854368932038835
$ws.Range("B24").Value = 0.8387096774193549
$ws.Range("C24").Value = 0.8387096774193549
$ws.Range("D24").Value = 0.8387096774193549
$ws.Range("E24").Value = 0.8387096774193549
$ws.Range("B25").Value = 0.836111111111111
$ws.Range("C25").Value = 0.8473684210526315
$ws.Range("D25").Value = 0.836823020236285
$ws.Range("B26").Value = 0.8508363201911587
$ws.Range("C26").Value = 0.8387096774193549
$ws.Range("D26").Value = 0.8400303374475039
